# This script updates the ObjTables header metadata strings found in cell A1
# (and A2 of the table-of-contents sheet) of every worksheet, matching the
# obj_tables / wc_lang library upgrade: attribute names are lower-cased
# (Type->type, Id->id, ObjTablesVersion->objTablesVersion) and a schema=
# and tableFormat= attribute are appended.
$wb = $excel.ActiveWorkbook

# Table of contents sheet: update the !!!ObjTables header (A1) and the
# !!ObjTables schema header (A2)
$toc = $wb.Worksheets.Item("!!_Table of contents")
$toc.Cells.Item(1,1).Value = "!!!ObjTables objTablesVersion='0.0.8'"
$toc.Cells.Item(2,1).Value = "!!ObjTables type='Schema' objTablesVersion='0.0.8' tableFormat='row'"

# Data sheets: update the !!ObjTables header in cell A1
$ws = $wb.Worksheets.Item("!!Model")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='Model' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='column'"
$ws = $wb.Worksheets.Item("!!Taxon")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='Taxon' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='column'"
$ws = $wb.Worksheets.Item("!!Environment")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='Environment' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='column'"
$ws = $wb.Worksheets.Item("!!Submodels")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='Submodel' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'"
$ws = $wb.Worksheets.Item("!!Compartments")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='Compartment' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'"
$ws = $wb.Worksheets.Item("!!Species types")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='SpeciesType' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'"
$ws = $wb.Worksheets.Item("!!Species")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='Species' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'"
$ws = $wb.Worksheets.Item("!!Init species concentrations")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='DistributionInitConcentration' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'"
$ws = $wb.Worksheets.Item("!!Observables")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='Observable' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'"
$ws = $wb.Worksheets.Item("!!Functions")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='Function' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'"
$ws = $wb.Worksheets.Item("!!Reactions")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='Reaction' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'"
$ws = $wb.Worksheets.Item("!!Rate laws")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='RateLaw' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'"
$ws = $wb.Worksheets.Item("!!dFBA objectives")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='DfbaObjective' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'"
$ws = $wb.Worksheets.Item("!!dFBA objective reactions")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='DfbaObjReaction' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'"
$ws = $wb.Worksheets.Item("!!dFBA objective species")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='DfbaObjSpecies' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'"
$ws = $wb.Worksheets.Item("!!Parameters")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='Parameter' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'"
$ws = $wb.Worksheets.Item("!!Stop conditions")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='StopCondition' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'"
$ws = $wb.Worksheets.Item("!!Observations")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='Observation' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'"
$ws = $wb.Worksheets.Item("!!Observation sets")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='ObservationSet' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'"
$ws = $wb.Worksheets.Item("!!Conclusions")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='Conclusion' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'"
$ws = $wb.Worksheets.Item("!!References")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='Reference' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'"
$ws = $wb.Worksheets.Item("!!Authors")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='Author' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'"
$ws = $wb.Worksheets.Item("!!Changes")
$ws.Cells.Item(1,1).Value = "!!ObjTables type='Data' id='Change' objTablesVersion='0.0.8' schema='wc_lang' tableFormat='row'"
